$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark (it currently sits at the end of the
#    "Rủi ro chức năng..." paragraph). Word will re-create it at the location of
#    the next edit, so we delete it here and re-add it where the new edit happens.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Delete the whole "+) Giao diện quản trị: ..." paragraph (its text and its
#    paragraph mark), merging it into the following paragraph
#    ("+) Đa người dùng và đầu vào đồng thời: ...").
$delRange = $d.Content
$delRange.Find.Execute(
    "+) Giao diện quản trị: Xây dựng một giao diện quản trị cho quản trị viên, cho phép quản trị viên thêm cấp độ truy cập mới, thay đổi mô tả phân loại và thực hiện các sửa đổi khác trong hệ thống báo cáo hoặc kho dữ liệu tập trung.^p",
    $true, $false, $false, $false, $false, $false, 1, $false, "", 2) | Out-Null

# 3) Drop a new "_GoBack" bookmark exactly where the edit happened: at the very
#    start of the now-adjacent "+) Đa người dùng..." paragraph.
$markRange = $d.Content
$markRange.Find.Execute(
    "+) Đa người dùng", $true, $false, $false, $false, $false, $false, 1, $false, "", 0) | Out-Null
$goBack = $d.Range($markRange.Start, $markRange.Start)
$d.Bookmarks.Add("_GoBack", $goBack) | Out-Null
